$d = $word.ActiveDocument
Write-Output "Paragraphs: $($d.Paragraphs.Count)"
